$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the timestamp (column O) for every data row (rows 2 through 405)
# to the new scrape time.
for ($r = 2; $r -le 405; $r++) {
    $ws.Cells.Item($r, 15).Value2 = "2022-12-22 12:56:18"
}

# Update the productAriaLabel text for row 398 (Leisi Blätterteig) to mention
# that it is out of stock online.
$ws.Cells.Item(398, 13).Value2 = "Leisi Blätterteig glutenfrei rund ausgewallt Ø32cm - Online kein Bestand 20% ab 2 Aktion 4.95 Schweizer Franken"

# Remove the "Biber mit Sujet 2x90g" row (row 404); the following row
# ("Bonne Maman Madeleine 7 Stück") shifts up to become the new row 404,
# and the used range shrinks from A1:O405 to A1:O404.
$ws.Rows(404).Delete()
